# "Hopeless calibration of Kazakhstan"
#
# Updates a few calibration constants on the "constants" sheet, drops the
# separator border that used to sit above the "susceptible_fully" row,
# adds a brand-new "age_breakpoints" parameter row, and removes the
# (now stale) data-validation rules on the "time_variants" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "constants" sheet: recalibrate a few values
# ---------------------------------------------------------------------
$constants = $wb.Worksheets.Item("constants")
$constants.Activate()

# tb_n_contact (point estimate) 11 -> 25
$constants.Range("B2").Value = 25

# program_prop_death_reporting (point estimate) 0.7 -> 1
$constants.Range("B3").Value = 1

# susceptible_fully (point estimate) 11,000,000 -> 15,000,000
$constants.Range("B4").Value = 15000000

# The row for susceptible_fully used to be set off from the row above by a
# thin top border; that separator line goes away so it matches the row
# below it (active).
$constants.Range("A4:E4").Borders.Item(8).LineStyle = -4142

# Add the new "age_breakpoints" parameter row, carrying over the same
# formatting used by the other rows in this block.
$constants.Range("A5").Copy()
$constants.Range("A6:E6").PasteSpecial(-4122)
$constants.Range("A6").Value = "age_breakpoints"
$constants.Range("B6").Value = 5
$constants.Range("C6").Value = 15

$constants.Range("A4:XFD6").Select()

# ---------------------------------------------------------------------
# 2) "time_variants" sheet: drop the old data-validation rules and leave
#    it as the active sheet/selection, matching where editing finished.
# ---------------------------------------------------------------------
$timeVariants = $wb.Worksheets.Item("time_variants")
$timeVariants.Activate()
$timeVariants.Cells.Validation.Delete()
$timeVariants.Range("A2:XFD2").Select()
